# Append two empty paragraphs (matching the "Press Enter twice" pilcrow
# formatting Word leaves behind) to the very end of the document body,
# right before the sectPr, each carrying the same run-formatting mark
# (Times New Roman, 12pt) as the rest of the document.

$d = $word.ActiveDocument

$blankParagraphXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
  "<w:pPr>" +
    "<w:rPr>" +
      "<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>" +
      "<w:sz w:val='24'/>" +
      "<w:szCs w:val='24'/>" +
    "</w:rPr>" +
  "</w:pPr>" +
"</w:p>"

# Collapse to the end of the story and insert the blank-paragraph XML
# twice so two new (empty) paragraphs land right before the section
# break, after the existing final paragraph.
$end1 = $d.Content
$end1.Collapse(0)
$null = $end1.InsertXML($blankParagraphXml)

$end2 = $d.Content
$end2.Collapse(0)
$null = $end2.InsertXML($blankParagraphXml)
